$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "datetimeFigureOut" field: 31/10/2022 -> 07/11/2022
#    This cached field text lives on the slide master and on every one of its
#    custom (slide) layouts. Walk them all and update the matching shape.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "31/10/2022") {
            $shp.TextFrame.TextRange.Text = "07/11/2022"
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "31/10/2022") {
                $shp.TextFrame.TextRange.Text = "07/11/2022"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) EDU-CIAA-NXP pinout table (slide 1, "85 Tabla"): re-assign several of the
#    GPIO / LCD labels in the middle column.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        if ($shp.Name -eq "85 Tabla") {
            $tableShape = $shp
        }
    }
}

$tbl = $tableShape.Table

$tbl.Cell(2,2).Shape.TextFrame.TextRange.Text  = "GPIO5"
$tbl.Cell(4,2).Shape.TextFrame.TextRange.Text  = "GPIO7"
$tbl.Cell(6,2).Shape.TextFrame.TextRange.Text  = "GPIO1"
$tbl.Cell(7,2).Shape.TextFrame.TextRange.Text  = "GPIO8"
$tbl.Cell(8,2).Shape.TextFrame.TextRange.Text  = "TCOL0"
$tbl.Cell(10,2).Shape.TextFrame.TextRange.Text = "GPIO4"
$tbl.Cell(11,2).Shape.TextFrame.TextRange.Text = "LCD2"
